$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted at row 142 ("Mercado Mayorista Lo Valledor
# de Santiago - Arveja Verde"). This pushes every existing record on or after
# row 142 down by one row (old row N -> new row N+1), which is exactly what
# Excel's native row-insert does, so the rest of the sheet needs no further
# editing. Only the brand-new top row (142) needs its values filled in.

$ws.Rows.Item(142).Insert()

$ws.Cells.Item(142, 1).Value = 6
$ws.Cells.Item(142, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(142, 3).Value = "Metropolitana"
$ws.Cells.Item(142, 4).Value = 44719
$ws.Cells.Item(142, 5).Value = 13
$ws.Cells.Item(142, 6).Value = 100112022
$ws.Cells.Item(142, 7).Value = "Arveja Verde"
$ws.Cells.Item(142, 8).Value = "Sin especificar"
$ws.Cells.Item(142, 9).Value = "Primera"
$ws.Cells.Item(142, 10).Value = 380
$ws.Cells.Item(142, 11).Value = 35000
$ws.Cells.Item(142, 12).Value = 37000
$ws.Cells.Item(142, 13).Value = 36211
$ws.Cells.Item(142, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(142, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(142, 16).Value = 1448
$ws.Cells.Item(142, 17).Value = 25
$ws.Cells.Item(142, 18).Value = "Hortaliza"
